# Referenties nieuwe CERA-bonnen (bis)
# - Remove sheets "08935" and "08936" (old voucher references no longer used)
# - Add a new sheet "08955" at the end, cloned from the last remaining
#   template sheet so it keeps the same layout/format/autofilter
# - Make the new sheet the active tab (it was the last tab that was selected)

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# Remove the two retired voucher-reference sheets
$wb.Worksheets.Item("08935").Delete()
$wb.Worksheets.Item("08936").Delete()

# Clone the last remaining template sheet ("08954") to create the new
# "08955" sheet with identical formatting / autofilter / frozen pane
$template = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "08955"

# Recreate the hidden _FilterDatabase defined name that Excel keeps per
# autoFiltered sheet, matching the other sheets in the workbook
$newSheet.Names.Add("_xlnm._FilterDatabase", "='08955'!`$A`$1:`$B`$7")

# The newest sheet becomes the active / selected tab
$newSheet.Activate()
$newSheet.Select()
